$d = $word.ActiveDocument

$replacements = @(
    @("804÷3=", "926÷6="),
    @("232÷4=", "277÷8="),
    @("473÷9=", "899÷2="),
    @("482÷8=", "549÷5="),
    @("400÷5=", "597÷2="),
    @("225÷8=", "657÷2="),
    @("974÷4=", "990÷9="),
    @("612÷2=", "579÷2="),
    @("562÷3=", "640÷6="),
    @("900÷5=", "668÷9="),
    @("273÷7=", "594÷7="),
    @("485÷9=", "702÷8="),
    @("956÷8=", "856÷8="),
    @("581÷5=", "655÷3="),
    @("403÷8=", "750÷9="),
    @("995÷7=", "390÷6="),
    @("968÷9=", "799÷7="),
    @("921÷7=", "708÷3="),
    @("457÷8=", "569÷9="),
    @("820÷2=", "221÷9="),
    @("420÷6=", "330÷2="),
    @("214÷9=", "110÷7="),
    @("965÷9=", "747÷2="),
    @("679÷2=", "318÷7="),
    @("499÷7=", "624÷5="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying $($replacements.Count) replacements"
